$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# --- Text change: label in C10 gets a longer description ---
$ws.Cells.Item(10, 3).Value = "Cost Savings from Violation Change from 2019-2020"

# --- Column width changes (A wider, C wider) ---
# ColumnWidth is quantized internally to the nearest pixel, so these inputs are
# chosen to land as close as possible to the desired stored widths.
$ws.Columns.Item(1).ColumnWidth = 39.917
$ws.Columns.Item(3).ColumnWidth = 31.417

# --- Wrap text for column A data cells (A1:A7) and the C10 label ---
$ws.Range("A1:A7").WrapText = $true
$ws.Range("C10").WrapText = $true

# --- Row heights: taller header/data rows, and a taller note row at the bottom ---
$ws.Range("A1:A7").RowHeight = 18
$ws.Rows.Item(10).RowHeight = 30

# --- Selection moves to A10 ---
$ws.Range("A10").Select()
